$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Header text updates (shared strings used by A8 "Volume 31 Number 6"
# and C9 "Report Covering the Week 2/5/2024 Through 2/11/2024")
# -----------------------------------------------------------------------

# A8: "Volume 31   Number  6" -> "...7"  (only the trailing issue number changes)
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$volLen = $volText.Length
$volCell.Characters($volLen, 1).Text = "7"

# C9: "Report Covering the Week  2/5/2024  Through  2/11/2024"
#     -> "...2/12/2024  Through  2/18/2024"
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value2
$pos1 = $weekText.IndexOf("2/5/2024") + 1
$weekCell.Characters($pos1, 8).Text = "2/12/2024"

$weekText2 = $weekCell.Value2
$pos2 = $weekText2.IndexOf("2/11/2024") + 1
$weekCell.Characters($pos2, 9).Text = "2/18/2024"

# -----------------------------------------------------------------------
# Weekly crime-statistics table (rows 14-27): new week's figures.
# Some cells flip between a numeric value and the sheet's placeholder
# shared strings ("0" / "***.*") used for zero-base / undefined %-change
# cells, so those cells first copy formatting+type from an untouched
# donor cell elsewhere in the same table (row 28) that already carries
# the right style, then (for numeric targets) get their real value
# written on top.
# -----------------------------------------------------------------------

# Row 14
$ws.Range("C28").Copy($ws.Range("D14"))
$ws.Range("E28").Copy($ws.Range("E14"))
$ws.Range("G14").Value2 = 1

# Row 15
$ws.Range("J28").Copy($ws.Range("D15"))
$ws.Range("D15").Value2 = 1
$ws.Range("K28").Copy($ws.Range("E15"))
$ws.Range("E15").Value2 = -100
$ws.Range("J28").Copy($ws.Range("G15"))
$ws.Range("G15").Value2 = 1
$ws.Range("K28").Copy($ws.Range("H15"))
$ws.Range("H15").Value2 = 0
$ws.Range("J15").Value2 = 2
$ws.Range("K15").Value2 = 50
$ws.Range("K28").Copy($ws.Range("M15"))
$ws.Range("M15").Value2 = 200
$ws.Range("N15").Value2 = -25

# Row 16
$ws.Range("C16").Value2 = 8
$ws.Range("J28").Copy($ws.Range("D16"))
$ws.Range("D16").Value2 = 2
$ws.Range("K28").Copy($ws.Range("E16"))
$ws.Range("E16").Value2 = 300
$ws.Range("F16").Value2 = 16
$ws.Range("G16").Value2 = 8
$ws.Range("H16").Value2 = 100
$ws.Range("I16").Value2 = 23
$ws.Range("J16").Value2 = 20
$ws.Range("K16").Value2 = 15
$ws.Range("L16").Value2 = -4.166666666666
$ws.Range("M16").Value2 = -20.689655172413
$ws.Range("N16").Value2 = -75.531914893617

# Row 17
$ws.Range("C17").Value2 = 1
$ws.Range("D17").Value2 = 5
$ws.Range("E17").Value2 = -80
$ws.Range("G17").Value2 = 9
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 17
$ws.Range("J17").Value2 = 18
$ws.Range("K17").Value2 = -5.555555555555
$ws.Range("L17").Value2 = -19.047619047619
$ws.Range("M17").Value2 = 0
$ws.Range("N17").Value2 = -66.666666666666

# Row 18
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 1
$ws.Range("E18").Value2 = 200
$ws.Range("F18").Value2 = 10
$ws.Range("H18").Value2 = 233.333333333333
$ws.Range("I18").Value2 = 18
$ws.Range("J18").Value2 = 13
$ws.Range("K18").Value2 = 38.461538461538
$ws.Range("L18").Value2 = -43.75
$ws.Range("M18").Value2 = -14.285714285714
$ws.Range("N18").Value2 = -88.75

# Row 19
$ws.Range("C19").Value2 = 9
$ws.Range("D19").Value2 = 5
$ws.Range("E19").Value2 = 80
$ws.Range("F19").Value2 = 34
$ws.Range("G19").Value2 = 36
$ws.Range("H19").Value2 = -5.555555555555
$ws.Range("I19").Value2 = 55
$ws.Range("J19").Value2 = 74
$ws.Range("K19").Value2 = -25.675675675675
$ws.Range("L19").Value2 = -20.289855072463
$ws.Range("M19").Value2 = -15.384615384615
$ws.Range("N19").Value2 = -54.918032786885

# Row 20
$ws.Range("D20").Value2 = 1
$ws.Range("E20").Value2 = 100
$ws.Range("F20").Value2 = 7
$ws.Range("G20").Value2 = 12
$ws.Range("H20").Value2 = -41.666666666666
$ws.Range("I20").Value2 = 8
$ws.Range("J20").Value2 = 18
$ws.Range("K20").Value2 = -55.555555555555
$ws.Range("L20").Value2 = 33.333333333333
$ws.Range("M20").Value2 = 166.666666666667
$ws.Range("N20").Value2 = -93.495934959349

# Row 21
$ws.Range("C21").Value2 = 23
$ws.Range("D21").Value2 = 15
$ws.Range("E21").Value2 = 53.333333333333
$ws.Range("F21").Value2 = 77
$ws.Range("G21").Value2 = 70
$ws.Range("H21").Value2 = 10
$ws.Range("I21").Value2 = 124
$ws.Range("J21").Value2 = 147
$ws.Range("K21").Value2 = -15.646258503401
$ws.Range("L21").Value2 = -18.954248366013
$ws.Range("M21").Value2 = -9.489051094890
$ws.Range("N21").Value2 = -77.777777777777

# Row 22
$ws.Range("J28").Copy($ws.Range("C22"))
$ws.Range("C22").Value2 = 1
$ws.Range("F22").Value2 = 4
$ws.Range("G22").Value2 = 2
$ws.Range("H22").Value2 = 100
$ws.Range("I22").Value2 = 6
$ws.Range("K22").Value2 = 100
$ws.Range("M22").Value2 = 100

# Row 23
$ws.Range("C23").Value2 = 1
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = -66.666666666666
$ws.Range("F23").Value2 = 8
$ws.Range("G23").Value2 = 11
$ws.Range("H23").Value2 = -27.272727272727
$ws.Range("I23").Value2 = 12
$ws.Range("J23").Value2 = 15
$ws.Range("K23").Value2 = -20
$ws.Range("M23").Value2 = -14.285714285714

# Row 24
$ws.Range("C24").Value2 = 17
$ws.Range("D24").Value2 = 30
$ws.Range("E24").Value2 = -43.333333333333
$ws.Range("F24").Value2 = 95
$ws.Range("G24").Value2 = 144
$ws.Range("H24").Value2 = -34.027777777777
$ws.Range("I24").Value2 = 177
$ws.Range("J24").Value2 = 295
$ws.Range("K24").Value2 = -40
$ws.Range("L24").Value2 = -29.482071713147
$ws.Range("M24").Value2 = 34.090909090909

# Row 25
$ws.Range("C25").Value2 = 3
$ws.Range("D25").Value2 = 7
$ws.Range("E25").Value2 = -57.142857142857
$ws.Range("F25").Value2 = 22
$ws.Range("G25").Value2 = 21
$ws.Range("H25").Value2 = 4.761904761904
$ws.Range("I25").Value2 = 40
$ws.Range("J25").Value2 = 39
$ws.Range("K25").Value2 = 2.564102564102
$ws.Range("L25").Value2 = 2.564102564102
$ws.Range("M25").Value2 = 0

# Row 26
$ws.Range("G26").Value2 = 2
$ws.Range("H26").Value2 = -50
$ws.Range("J26").Value2 = 3
$ws.Range("K26").Value2 = 0

# Row 27
$ws.Range("C28").Copy($ws.Range("C27"))
$ws.Range("C28").Copy($ws.Range("D27"))
$ws.Range("E28").Copy($ws.Range("E27"))
$ws.Range("F27").Value2 = 3
$ws.Range("G27").Value2 = 3
